$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.058.05"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.886.41"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7370"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.23"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3170"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07200"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.78"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08332"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7573"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.403"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.888.09"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.05"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.159"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.045.38"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "249.82"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.58"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007859"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.147.16"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.896"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9994"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1565"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.285"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.71"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.051"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.477"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.577"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.535"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.204"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05338"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.251"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7707"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9988"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.720"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01964"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.758"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4588"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.040"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.088.41"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.39"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8738"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.65"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.859"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.586"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.544"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.037.81"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
